# files generation fixes & upd
#
# Row 2 used to hold a single placeholder cell (B2 = "default", rendered with
# the bold/bordered "header" look re-used from a leftover style). The sheet
# is regenerated so that row 2 and a brand-new row 3 are fully populated
# (A:O) with plain, unstyled text values: the row number repeated in most
# columns and a single space in columns C, G, H and I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 carried the one-off thin-border style (style index 2). Strip its border
# so it collapses back onto the sheet's default (unstyled) cell format,
# matching every other cell we are about to write.
$ws.Range("B2").Borders.LineStyle = 0

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O")
$spaceColumns = @("C","G","H","I")

# A scratch cell well outside the printed area: used to coerce a numeric-
# looking string ("2", "3") into a genuine text value via Copy/PasteSpecial
# (values only), so the cell keeps the default style and isn't reinterpreted
# as a number.
$scratch = $ws.Range("Z1")

foreach ($r in 2,3) {
    foreach ($c in $columns) {
        if ($spaceColumns -contains $c) {
            $scratch.Formula = '=" "'
        } else {
            $scratch.Formula = "=""$r"""
        }
        $scratch.Copy()
        $ws.Range("$c$r").PasteSpecial(-4163)
    }
}

$scratch.Clear()

Write-Host "done"
